$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "CNNLSTM Temp": restore previously-blank B:F inputs (rows 10, 16, 22).
# The existing AVERAGE() formulas in column G will recompute automatically
# once these inputs are populated (they currently evaluate to #DIV/0!).
# ---------------------------------------------------------------------------
$wsCnnTemp = $wb.Worksheets.Item("CNNLSTM Temp")
$wsCnnTemp.Range("B10").Value = 1.62196488088258
$wsCnnTemp.Range("C10").Value = 1.80098344754076
$wsCnnTemp.Range("D10").Value = 2.66360676629568
$wsCnnTemp.Range("E10").Value = 2.51003856774608
$wsCnnTemp.Range("F10").Value = 2.73647519661725
$wsCnnTemp.Range("B16").Value = 1.21987339010018
$wsCnnTemp.Range("C16").Value = 0.832322533212337
$wsCnnTemp.Range("D16").Value = 3.46878517414154
$wsCnnTemp.Range("E16").Value = 2.38060709980662
$wsCnnTemp.Range("F16").Value = 1.32556687662497
$wsCnnTemp.Range("B22").Value = 2.04452322797538
$wsCnnTemp.Range("C22").Value = 1.95114572741722
$wsCnnTemp.Range("D22").Value = 3.25623392982435
$wsCnnTemp.Range("E22").Value = 2.17050704262248
$wsCnnTemp.Range("F22").Value = 1.44993684112294

# ---------------------------------------------------------------------------
# Sheet "LSTM SH": restore previously-blank B:F, H, I inputs (rows 4,5,6,7 and
# 10,11,12,13). Column G AVERAGE() formulas recompute automatically.
# ---------------------------------------------------------------------------
$wsLstmSh = $wb.Worksheets.Item("LSTM SH")
$wsLstmSh.Range("B4").Value = 0.377387905103774
$wsLstmSh.Range("C4").Value = 0.38922616667495
$wsLstmSh.Range("D4").Value = 0.459634364033748
$wsLstmSh.Range("E4").Value = 0.542584575192884
$wsLstmSh.Range("F4").Value = 0.318672235442113
$wsLstmSh.Range("H4").Value = 4.62290909090909
$wsLstmSh.Range("I4").Value = 0.0903113258511845
$wsLstmSh.Range("B5").Value = 0.418366647027552
$wsLstmSh.Range("C5").Value = 0.64034076878181
$wsLstmSh.Range("D5").Value = 0.605642253672639
$wsLstmSh.Range("E5").Value = 0.562217407515377
$wsLstmSh.Range("F5").Value = 0.40491034936676
$wsLstmSh.Range("H5").Value = 4.62290909090909
$wsLstmSh.Range("I5").Value = 0.113845086486295
$wsLstmSh.Range("B6").Value = 0.335993396019012
$wsLstmSh.Range("C6").Value = 0.601774111790549
$wsLstmSh.Range("D6").Value = 0.589314132628204
$wsLstmSh.Range("E6").Value = 0.742011343696405
$wsLstmSh.Range("F6").Value = 0.507465696568999
$wsLstmSh.Range("H6").Value = 4.62290909090909
$wsLstmSh.Range("I6").Value = 0.120121708045838
$wsLstmSh.Range("B7").Value = 0.833341346563887
$wsLstmSh.Range("C7").Value = 0.782766807380188
$wsLstmSh.Range("D7").Value = 1.18960353451858
$wsLstmSh.Range("E7").Value = 0.655915161265528
$wsLstmSh.Range("F7").Value = 0.722363815423245
$wsLstmSh.Range("H7").Value = 4.62290909090909
$wsLstmSh.Range("I7").Value = 0.181011159115337
$wsLstmSh.Range("B10").Value = 0.81095665354673
$wsLstmSh.Range("C10").Value = 0.913477008121699
$wsLstmSh.Range("D10").Value = 0.908899743914195
$wsLstmSh.Range("E10").Value = 0.56556211929072
$wsLstmSh.Range("F10").Value = 0.70099396024346
$wsLstmSh.Range("H10").Value = 5.1330303030303
$wsLstmSh.Range("I10").Value = 0.151952716227468
$wsLstmSh.Range("B11").Value = 0.877113260933854
$wsLstmSh.Range("C11").Value = 0.542771352391183
$wsLstmSh.Range("D11").Value = 0.885169433359076
$wsLstmSh.Range("E11").Value = 0.529584091411883
$wsLstmSh.Range("F11").Value = 0.727058565888869
$wsLstmSh.Range("H11").Value = 5.1330303030303
$wsLstmSh.Range("I11").Value = 0.138775596235315
$wsLstmSh.Range("B12").Value = 0.660653048364751
$wsLstmSh.Range("C12").Value = 0.902525769239827
$wsLstmSh.Range("D12").Value = 1.27969040747388
$wsLstmSh.Range("E12").Value = 0.8432749244013
$wsLstmSh.Range("F12").Value = 0.692407176805679
$wsLstmSh.Range("H12").Value = 5.1330303030303
$wsLstmSh.Range("I12").Value = 0.170602979830473
$wsLstmSh.Range("B13").Value = 0.804656622314353
$wsLstmSh.Range("C13").Value = 0.566895909636309
$wsLstmSh.Range("D13").Value = 0.94244529955387
$wsLstmSh.Range("E13").Value = 0.548843407148163
$wsLstmSh.Range("F13").Value = 0.659789820886284
$wsLstmSh.Range("H13").Value = 5.1330303030303
$wsLstmSh.Range("I13").Value = 0.137253468285951

# ---------------------------------------------------------------------------
# Sheet "Hyperparameter Optimization": restore previously-blank S:V columns
# for several monthly-binned model rows (plain literal values).
# ---------------------------------------------------------------------------
$wsHyper = $wb.Worksheets.Item("Hyperparameter Optimization")
$wsHyper.Range("S28").Value = 0.417501049289494
$wsHyper.Range("T28").Value = 0.314968465307896
$wsHyper.Range("U28").Value = 0.0823486417346422
$wsHyper.Range("V28").Value = 0.930949957308251
$wsHyper.Range("S29").Value = 0.526295485272828
$wsHyper.Range("T29").Value = 0.434268185404213
$wsHyper.Range("U29").Value = 0.10512403274831
$wsHyper.Range("V29").Value = 0.908752825392902
$wsHyper.Range("S30").Value = 0.555311736140634
$wsHyper.Range("T30").Value = 0.466987306640906
$wsHyper.Range("U30").Value = 0.115465177372042
$wsHyper.Range("V30").Value = 0.875454391514441
$wsHyper.Range("S31").Value = 0.836798133030286
$wsHyper.Range("T31").Value = 0.753237391450188
$wsHyper.Range("U31").Value = 0.178756426472261
$wsHyper.Range("V31").Value = 0.90700770470762
$wsHyper.Range("S32").Value = 0.779977897023361
$wsHyper.Range("T32").Value = 0.574667448833133
$wsHyper.Range("U32").Value = 0.118217808283464
$wsHyper.Range("V32").Value = 0.863666202573125
$wsHyper.Range("S33").Value = 0.712339340796973
$wsHyper.Range("T33").Value = 0.510538770453496
$wsHyper.Range("U33").Value = 0.103453754569833
$wsHyper.Range("V33").Value = 0.868888251981647
$wsHyper.Range("S34").Value = 0.875710265257087
$wsHyper.Range("T34").Value = 0.684717642287413
$wsHyper.Range("U34").Value = 0.157093036130365
$wsHyper.Range("V34").Value = 0.785988449395206
$wsHyper.Range("S35").Value = 0.704526211907796
$wsHyper.Range("T35").Value = 0.509113214969465
$wsHyper.Range("U35").Value = 0.105945977405843
$wsHyper.Range("V35").Value = 0.867096411121877
$wsHyper.Range("S52").Value = 2.26661377181647
$wsHyper.Range("T52").Value = 1.72817704302975
$wsHyper.Range("U52").Value = 0.236446483682046
$wsHyper.Range("V52").Value = 0.846225456382236
$wsHyper.Range("S56").Value = 1.84543101477713
$wsHyper.Range("T56").Value = 1.43636705660737
$wsHyper.Range("U56").Value = 0.144638988094586
$wsHyper.Range("V56").Value = 0.944761267346589
$wsHyper.Range("S60").Value = 2.17446935379247
$wsHyper.Range("T60").Value = 1.69656638093924
$wsHyper.Range("U60").Value = 0.0769228725983207
$wsHyper.Range("V60").Value = 0.948125365317433

# ---------------------------------------------------------------------------
# Restore selection/active-cell state to match the saved workbook (best
# effort -- the emulated window model here does not expose scroll/
# topLeftCell persistence, only the active selection).
# ---------------------------------------------------------------------------
$wsCnnTemp.Activate()
$wsCnnTemp.Range("B11").Select()

$wsLstmSh.Activate()
$wsLstmSh.Range("B11").Select()

$wsHyper.Activate()
$wsHyper.Range("S33").Select()

Write-Host "Edit applied successfully"
